# Junction_Flooding_83.xlsx edit:
#  - Row 5 data values recomputed with "custom accuracy" (≈2 decimal places)
#  - Row 6 (the extra data row) removed
#  - Column T (20th column) width trimmed from 9 to 8 characters

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 5 with the new, lower-precision readings ---
$ws.Cells.Item(5, 1).Value = 44781.9027662037
$ws.Cells.Item(5, 2).Value = 17.3
$ws.Cells.Item(5, 3).Value = 12.65
$ws.Cells.Item(5, 4).Value = 1.13
$ws.Cells.Item(5, 5).Value = 37.47
$ws.Cells.Item(5, 6).Value = 30.75
$ws.Cells.Item(5, 7).Value = 13.61
$ws.Cells.Item(5, 8).Value = 49.02
$ws.Cells.Item(5, 9).Value = 20.94
$ws.Cells.Item(5, 10).Value = 9.24
$ws.Cells.Item(5, 11).Value = 13.74
$ws.Cells.Item(5, 12).Value = 15.08
$ws.Cells.Item(5, 13).Value = 15.83
$ws.Cells.Item(5, 14).Value = 4.35
$ws.Cells.Item(5, 15).Value = 13.53
$ws.Cells.Item(5, 16).Value = 19.2
$ws.Cells.Item(5, 17).Value = 11.46
$ws.Cells.Item(5, 18).Value = 0.83
$ws.Cells.Item(5, 19).Value = 0.74
$ws.Cells.Item(5, 20).Value = 198.67
$ws.Cells.Item(5, 21).Value = 37.74
$ws.Cells.Item(5, 22).Value = 12.49
$ws.Cells.Item(5, 23).Value = 25.29
$ws.Cells.Item(5, 24).Value = 13.36
$ws.Cells.Item(5, 25).Value = 2.06
$ws.Cells.Item(5, 26).Value = 24.18
$ws.Cells.Item(5, 27).Value = 11.03
$ws.Cells.Item(5, 28).Value = 9.83
$ws.Cells.Item(5, 29).Value = 11.55
$ws.Cells.Item(5, 30).Value = 15.72
$ws.Cells.Item(5, 31).Value = 0.5600000000000001
$ws.Cells.Item(5, 32).Value = 44.23
$ws.Cells.Item(5, 33).Value = 7.03
$ws.Cells.Item(5, 34).Value = 15.62

# --- Remove the now-obsolete row 6 entirely (shifts dimension to A1:AH5) ---
$ws.Cells.Item(6, 1).EntireRow.Delete()

# --- Narrow column T (20th column) from width 9 down to width 8 ---
# Excel pads the COM "ColumnWidth" by 5/6 of a character when storing the
# raw OOXML width, so subtract that padding to land on an exact width of 8.
$ws.Columns.Item(20).ColumnWidth = 8 - 5/6
